$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 229, shifting existing rows 229:267 down to 230:268.
$ws.Rows.Item(229).Insert()

# Populate the newly inserted row 229 with this week's data. Metadata columns
# (A,B,C,E,F,G,H,I,N,O,Q,R) mirror the row immediately below (the former row
# 229, now row 230); only the date and the price columns change.
$ws.Range("A229").Value = 8
$ws.Range("B229").Value = "Terminal La Palmera de La Serena"
$ws.Range("C229").Value = "Coquimbo"
$ws.Range("D229").Value = 45258
$ws.Range("D229").NumberFormat = $ws.Range("D230").NumberFormat
$ws.Range("E229").Value = 4
$ws.Range("F229").Value = 100112044
$ws.Range("G229").Value = "Perejil"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 2000
$ws.Range("K229").Value = 2300
$ws.Range("L229").Value = 2500
$ws.Range("M229").Value = 2400
$ws.Range("N229").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O229").Value = "Provincia del Elquí"
$ws.Range("P229").Value = 1600
$ws.Range("Q229").Value = 1.5
$ws.Range("R229").Value = "Hortaliza"
